$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '56.846.20'
$ws.Range("E2").Value = '  -3.70%  '

# Row 3
$ws.Range("D3").Value = '2.532.12'
$ws.Range("E3").Value = '  -5.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '514.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.76%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.79%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.557'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.32%  '

# Row 9
$ws.Range("D9").Value = '2.535.83'
$ws.Range("E9").Value = '  -5.13%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.19%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0988'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.99%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.321'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.24%  '

# Row 13
$ws.Range("E13").Value = '  -0.41%  '

# Row 14
$ws.Range("D14").Value = '2.977.52'
$ws.Range("E14").Value = '  -4.92%  '

# Row 15
$ws.Range("D15").Value = '56.841.03'
$ws.Range("E15").Value = '  -3.68%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.21%  '

# Row 17
$ws.Range("E17").Value = '  -3.39%  '

# Row 18
$ws.Range("D18").Value = '2.537.52'
$ws.Range("E18").Value = '  -5.19%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '329.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.93%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.70%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.46%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.32%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.398'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.19%  '

# Row 28
$ws.Range("D28").Value = '2.651.21'
$ws.Range("E28").Value = '  -4.81%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.12%  '

# Row 30
$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '

# Row 31
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0739'
$ws.Range("E31").Value = '  -7.87%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.80%  '

# Row 33
$ws.Range("E33").Value = '  -2.84%  '

# Row 34
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.66%  '

# Row 35
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.81%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.22%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.835'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.64%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.81%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.816'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.22%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.14%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.62%  '

# Row 44
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.64%  '

# Row 45
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0948'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.33%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '263.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.55%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.575'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.77%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.59%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0515'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.79%  '

# Row 50
$ws.Range("D50").Value = '1.947.42'
$ws.Range("E50").Value = '  -5.19%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0219'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.04%  '
